$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "papaya"
$ws.Range("B6").Value = "Salmon tacos with papaya salsa"
$ws.Range("C6").Value = 16

$ws.Range("A7").Value = "passion fruit"
$ws.Range("B7").Value = "Glen of passion"
$ws.Range("C7").Value = 8

$ws.Range("A8").Value = "plum"
$ws.Range("B8").Value = "Mini plum crostatas"
$ws.Range("C8").Value = 13

$ws.Range("B2").Value = "Apple and cinnamon choux fritters"
$ws.Range("C2").Value = 9
